# Pre-app Predetermination letter.docx edit script
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert three additional empty "KeyHeadDetails" paragraphs just
#    before the paragraph that holds the four tab characters (the row
#    right after "Our Ref: <Primary Reference Number>").
# ---------------------------------------------------------------------
$tabsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "`t`t`t`t`r") {
        $tabsPara = $p
        break
    }
}
if ($tabsPara -ne $null) {
    $tabsPara.Range.InsertParagraphBefore()
    $tabsPara.Range.InsertParagraphBefore()
    $tabsPara.Range.InsertParagraphBefore()
}

# ---------------------------------------------------------------------
# 2. "NATIONAL PLANNING POLICY FRAMEWORK 2021" -> "...2023"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("NATIONAL PLANNING POLICY FRAMEWORK 2021", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "NATIONAL PLANNING POLICY FRAMEWORK 2023", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Italicise the "<Proposal Description>" placeholder paragraph.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "<Proposal Description>`r") {
        $p.Range.Font.Italic = $true
        $p.Range.Font.ItalicBi = $true
        break
    }
}

# ---------------------------------------------------------------------
# 4. NPPF paragraph-number updates.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("NPPF paragraph 194 says", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "NPPF paragraph 200 says", 2) | Out-Null

$d.Content.Find.Execute("NPPF paragraphs 199 - 202 place", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "NPPF paragraphs 205 - 208 place", 2) | Out-Null

$d.Content.Find.Execute("(NPPF paragraph 203)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(NPPF paragraph 209)", 2) | Out-Null

$d.Content.Find.Execute("NPPF paragraphs 190 and 197 and", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "NPPF paragraphs 195 and 203 and", 2) | Out-Null

$d.Content.Find.Execute("paragraph 205 of the NPPF", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "paragraph 211 of the NPPF", 2) | Out-Null

# ---------------------------------------------------------------------
# 5. Bold the second "<Casework Officer>" placeholder (the one in the
#    signature block, after "Yours sincerely").
# ---------------------------------------------------------------------
$seenFirst = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "<Casework Officer>`r") {
        if ($seenFirst) {
            $p.Range.Font.Bold = $true
            $p.Range.Font.BoldBi = $true
            break
        }
        $seenFirst = $true
    }
}

Write-Output "done"
